# Add season record columns (Wins, Losses, Ties) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting from an existing header cell (AC1) onto the
# new header cells so they match the bold/bordered/centered look of the
# other headers in row 1.
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player record in this sheet shares the same team season record.
for ($row = 2; $row -le 47; $row++) {
    $ws.Cells.Item($row, 30).Value = 71
    $ws.Cells.Item($row, 31).Value = 91
    $ws.Cells.Item($row, 32).Value = 0
}

Write-Output "done"
